# Scrum (Sprint Backlog) updated
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint Backlog")

# --- View state: scrolled/selected differently after the edit ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("E22").Select()

# --- Populate the "Actual" effort column (I) for rows 13-21, mirroring the
#     Planned effort column (H) that was already filled in ---
$ws.Range("I13").Value = 5
$ws.Range("I14").Value = 5
$ws.Range("I15").Value = 5
$ws.Range("I16").Value = 5
$ws.Range("I17").Value = 5
$ws.Range("I18").Value = 5
$ws.Range("I19").Value = 10
$ws.Range("I20").Value = 10
$ws.Range("I21").Value = 10

# --- Row 18 is repurposed from "DB Controller / DB Model" into a new
#     "Therapy Controller" task ---
$ws.Range("C18").Value = "Controller for therapy"
$ws.Range("D18").Value = "Implement therapy controller with state pattern"
$ws.Range("E18").Value = "Therapy Controller"
$ws.Range("H18").Value = 0
$ws.Range("L18").Value = "Therapy"

# --- The old "2.7" sprint-task number moves off row 19 (its task stays put,
#     but is no longer labelled here - it is relisted below as row 27) ---
$ws.Range("A19").ClearContents()

# --- Rename a duplicated/garbled user story description ---
$ws.Range("E21").Value = "Medication Model / Therapy Model"

# --- Sub-total row gains the "2.10" sprint marker and an "Actual" total ---
$ws.Range("A22").Value = "2.10"
$ws.Range("I22").Formula = "=SUM(I13:I21)"

$ws.Range("A23").ClearContents()

# --- Shift the trailing legend entries up one row, dropping the old row 25 ---
$ws.Range("L23").Value = "Add patient "
$ws.Range("L24").Value = "Add a medication"
$ws.Rows.Item(25).Delete()

# --- Re-add the "2.7 / Communication with Database" task, now two rows
#     below the sub-total, with the DB-controller wording typo fixed ---
$ws.Range("A27").Value = "2.7"
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = "Communication with Database"
$ws.Range("D27").Value = "Implement DB Controller and Model"
$ws.Range("E27").Value = "DB Controller / DB Model"
$ws.Range("F27").Value = "Stefan"
$ws.Range("G27").Value = "high"
$ws.Range("H27").Value = 5
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = "waiting"

$ws.Range("A2").Select()
